$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Bench1")

$ws1.Range("O3").Value = 5.5220000000000002
$ws1.Range("R3").Value = 8.0980000000000008
$ws1.Range("U3").Value = 8.5
$ws1.Range("O4").Value = 9.0920000000000005
$ws1.Range("R4").Value = 9.1660000000000004
$ws1.Range("U4").Value = 9.0359999999999996
$ws1.Range("O5").Value = 5.9989999999999997
$ws1.Range("R5").Value = 5.1230000000000002
$ws1.Range("U5").Value = 5.2809999999999997
$ws1.Range("O6").Value = 7.6769999999999996
$ws1.Range("R6").Value = 5.8109999999999999
$ws1.Range("U6").Value = 5.8890000000000002
$ws1.Range("O7").Value = 9.0359999999999996
$ws1.Range("R7").Value = 6.6139999999999999
$ws1.Range("U7").Value = 6.5419999999999998
$ws1.Range("O8").Value = 5.2539999999999996
$ws1.Range("R8").Value = 7.5640000000000001
$ws1.Range("U8").Value = 7.649
$ws1.Range("O9").Value = 5.8129999999999997
$ws1.Range("R9").Value = 8.3260000000000005
$ws1.Range("U9").Value = 8.6289999999999996
$ws1.Range("O10").Value = 6.5860000000000003
$ws1.Range("R10").Value = 9.5090000000000003
$ws1.Range("U10").Value = 9.4849999999999994
$ws1.Range("O11").Value = 7.351
$ws1.Range("R11").Value = 5.0060000000000002
$ws1.Range("U11").Value = 5.2359999999999998
$ws1.Range("O12").Value = 8.2550000000000008
$ws1.Range("R12").Value = 5.3250000000000002
$ws1.Range("U12").Value = 5.4669999999999996
$ws1.Range("O15").Value = 16.428999999999998
$ws1.Range("Q15").Value = 3276800000
$ws1.Range("R15").Value = 9.3260000000000005
$ws1.Range("U15").Value = 7.6970000000000001
$ws1.Range("O16").Value = 27.652000000000001
$ws1.Range("R16").Value = 6.2990000000000004
$ws1.Range("U16").Value = 7.4119999999999999
$ws1.Range("O17").Value = 38.869999999999997
$ws1.Range("R17").Value = 5.4420000000000002
$ws1.Range("T17").Value = 1638400000
$ws1.Range("U17").Value = 8.6519999999999992
$ws1.Range("O18").Value = 50.551000000000002
$ws1.Range("R18").Value = 7.7519999999999998
$ws1.Range("U18").Value = 5.9809999999999999
$ws1.Range("O19").Value = 61.862000000000002
$ws1.Range("R19").Value = 5.9269999999999996
$ws1.Range("U19").Value = 7.55
$ws1.Range("O20").Value = 73.462000000000003
$ws1.Range("R20").Value = 6.5949999999999998
$ws1.Range("U20").Value = 8.1010000000000009
$ws1.Range("O21").Value = 84.179000000000002
$ws1.Range("R21").Value = 7.883
$ws1.Range("U21").Value = 9.5730000000000004
$ws1.Range("O22").Value = 96.257000000000005
$ws1.Range("Q22").Value = 51200000
$ws1.Range("R22").Value = 9.6050000000000004
$ws1.Range("U22").Value = 7.3810000000000002
$ws1.Range("O23").Value = 107.12
$ws1.Range("R23").Value = 6.0549999999999997
$ws1.Range("U23").Value = 6.12
$ws1.Range("O24").Value = 118.057
$ws1.Range("R24").Value = 7.5369999999999999
$ws1.Range("U24").Value = 7.5839999999999996
$ws1.Range("O27").Value = 85.010999999999996
$ws1.Range("R27").Value = 6.2240000000000002
$ws1.Range("U27").Value = 5.6589999999999998
$ws1.Range("O28").Value = 95.278999999999996
$ws1.Range("R28").Value = 6.9470000000000001
$ws1.Range("U28").Value = 7.06
$ws1.Range("O29").Value = 105.724
$ws1.Range("Q29").Value = 51200000
$ws1.Range("R29").Value = 5.0709999999999997
$ws1.Range("U29").Value = 6.4649999999999999
$ws1.Range("O30").Value = 117.69799999999999
$ws1.Range("R30").Value = 7.0979999999999999
$ws1.Range("T30").Value = 102400000
$ws1.Range("U30").Value = 10.243
$ws1.Range("O31").Value = 127.82
$ws1.Range("R31").Value = 8.9440000000000008
$ws1.Range("U31").Value = 6.8550000000000004
$ws1.Range("O32").Value = 138.018
$ws1.Range("R32").Value = 5.1539999999999999
$ws1.Range("U32").Value = 8.75

# Update selection on Bench1 (was S26, now V26) and make Bench1 the active/selected sheet
$ws1.Activate() | Out-Null
$ws1.Range("V26").Select() | Out-Null
